$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: set column B values to 0 for rows 2,4,5,6,7
# and clear column C entirely for rows 2-7 (C column no longer populated) ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("B4").Value = 0
$wsRange.Range("B5").Value = 0
$wsRange.Range("B6").Value = 0
$wsRange.Range("B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- "Species qualification" sheet: Range Analysis (row5) species count set to 0 ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- "High Priority break-up" sheet: New High Species (no.) updated from 1 to 8 ---
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("D2").Value = 8
